$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the Price (D) column cells whose new values could
# otherwise be auto-converted to numbers by Excel, then restore default style so
# no stray formatting remains on the cells.

# Row 2
$ws.Range("D2").Value = '25.899.16'
$ws.Range("E2").Value = '  +0.52%  '

# Row 3
$ws.Range("D3").Value = '1.640.48'
$ws.Range("E3").Value = '  +0.99%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.50%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5084'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.35%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.21%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2598'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.76%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06469'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.82%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.95%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07834'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.85%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.667.64'
$ws.Range("E12").Value = '  +2.44%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.271'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.53%  '

# Row 14
$ws.Range("D14").Value = '1.864.62'
$ws.Range("E14").Value = '  +0.87%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5669'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.40%  '

# Row 16
$ws.Range("D16").Value = '0.0₅7707'

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.17%  '

# Row 18
$ws.Range("D18").Value = '25.897.88'
$ws.Range("E18").Value = '  +0.47%  '

# Row 19
$ws.Range("E19").Value = '  +0.26%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '195.01'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.63%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.403'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.15%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.46%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.213'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.01%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.003'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.16%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.765'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.36%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '138.26'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.70%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1240'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.15%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.875'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.18%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.34%  '

# Row 30
$ws.Range("E30").Value = '  +1.07%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05017'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.324'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.36%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.264'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.76%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.578'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.20%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.385'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.97%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9088'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.86%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.580'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.78%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5538'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.93%  '

# Row 39
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.131.67'
$ws.Range("E39").Value = '  -0.33%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01581'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.58%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9949'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.56%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.91'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.82%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.490'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8010'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.07%  '

# Row 45
$ws.Range("E45").Value = '  -2.53%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '55.80'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.10%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4238'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.06%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.734'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05047'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.53%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.9990'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.07%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.001'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
